# Updates recalculated profit-tracking figures (currentAveragePrice / Leve price / profit
# columns H:N) across several Leve rows on multiple job sheets, reflecting refreshed
# market-board data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 2935.8333
$ws.Range("I69").Value = 2920
$ws.Range("J69").Value = 3015
$ws.Range("K69").Value = 8760
$ws.Range("L69").Value = 9045
$ws.Range("M69").Value = -7886
$ws.Range("N69").Value = -10793
# Row 72
$ws.Range("H72").Value = 2935.8333
$ws.Range("I72").Value = 2920
$ws.Range("J72").Value = 3015
$ws.Range("K72").Value = 26280
$ws.Range("L72").Value = 27135
$ws.Range("M72").Value = -21912
$ws.Range("N72").Value = -35871
# Row 112
$ws.Range("H112").Value = 4156.6
$ws.Range("I112").Value = 970
$ws.Range("K112").Value = 2910
$ws.Range("M112").Value = -1802
# Row 129
$ws.Range("H129").Value = 1427.8
$ws.Range("J129").Value = 1852.5294
$ws.Range("L129").Value = 5557.5882
$ws.Range("N129").Value = -15557.5882
# Row 137
$ws.Range("H137").Value = 1210.6852
$ws.Range("I137").Value = 1189.0667
$ws.Range("K137").Value = 3567.2001
$ws.Range("M137").Value = -1017.2001

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 945.3333
$ws.Range("I74").Value = 849.51514
$ws.Range("J74").Value = 1296.6666
$ws.Range("K74").Value = 849.51514
$ws.Range("L74").Value = 1296.6666
$ws.Range("M74").Value = 24.48486000000003
$ws.Range("N74").Value = -3044.6666
# Row 77
$ws.Range("H77").Value = 945.3333
$ws.Range("I77").Value = 849.51514
$ws.Range("J77").Value = 1296.6666
$ws.Range("K77").Value = 4247.575699999999
$ws.Range("L77").Value = 6483.333000000001
$ws.Range("M77").Value = 120.4243000000006
$ws.Range("N77").Value = -15219.333

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2572.8667
$ws.Range("I134").Value = 2457.037
$ws.Range("K134").Value = 7371.110999999999
$ws.Range("M134").Value = -4836.110999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 261530.52
$ws.Range("I132").Value = 437241.53
$ws.Range("J132").Value = 2147.6191
$ws.Range("K132").Value = 1311724.59
$ws.Range("L132").Value = 6442.8573
$ws.Range("M132").Value = -1309194.59
$ws.Range("N132").Value = -11502.8573

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1724.9692
$ws.Range("I131").Value = 2840
$ws.Range("J131").Value = 1651.8524
$ws.Range("K131").Value = 8520
$ws.Range("L131").Value = 4955.5572
$ws.Range("M131").Value = -3480
$ws.Range("N131").Value = -15035.5572
# Row 132
$ws.Range("H132").Value = 1936.72
$ws.Range("J132").Value = 2208.2
$ws.Range("L132").Value = 19873.8
$ws.Range("N132").Value = -24933.8
# Row 138
$ws.Range("H138").Value = 2714.6428
$ws.Range("I138").Value = 916.5
$ws.Range("J138").Value = 4063.25
$ws.Range("K138").Value = 2749.5
$ws.Range("L138").Value = 12189.75
$ws.Range("M138").Value = 2390.5
$ws.Range("N138").Value = -22469.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2998.3333
$ws.Range("I80").Value = 2998.3333
$ws.Range("K80").Value = 2998.3333
$ws.Range("M80").Value = -2000.3333
# Row 83
$ws.Range("H83").Value = 2998.3333
$ws.Range("I83").Value = 2998.3333
$ws.Range("K83").Value = 14991.6665
$ws.Range("M83").Value = -9999.666499999999
# Row 122
$ws.Range("H122").Value = 2387.5454
$ws.Range("I122").Value = 1863.5238
$ws.Range("J122").Value = 3304.5833
$ws.Range("K122").Value = 5590.5714
$ws.Range("L122").Value = 9913.749899999999
$ws.Range("M122").Value = -3140.5714
$ws.Range("N122").Value = -14813.7499
# Row 132
$ws.Range("H132").Value = 2748.3333
$ws.Range("I132").Value = 1933.8334
$ws.Range("J132").Value = 4377.3335
$ws.Range("K132").Value = 5801.5002
$ws.Range("L132").Value = 13132.0005
$ws.Range("M132").Value = -3271.5002
$ws.Range("N132").Value = -18192.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1160.625
$ws.Range("I22").Value = 1168.6
$ws.Range("J22").Value = 1147.3334
$ws.Range("K22").Value = 1168.6
$ws.Range("L22").Value = 1147.3334
$ws.Range("M22").Value = -873.5999999999999
$ws.Range("N22").Value = -1737.3334
# Row 27
$ws.Range("H27").Value = 1160.625
$ws.Range("I27").Value = 1168.6
$ws.Range("J27").Value = 1147.3334
$ws.Range("K27").Value = 1168.6
$ws.Range("L27").Value = 1147.3334
$ws.Range("M27").Value = -1061.6
$ws.Range("N27").Value = -1361.3334
# Row 40
$ws.Range("H40").Value = 3103.4546
$ws.Range("I40").Value = 2913.8
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2913.8
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2777.8
$ws.Range("N40").Value = -5272
# Row 46
$ws.Range("H46").Value = 1633.3334
$ws.Range("I46").Value = 1533.3334
$ws.Range("K46").Value = 1533.3334
$ws.Range("M46").Value = -1345.3334
# Row 82
$ws.Range("H82").Value = 2666.6667
$ws.Range("J82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("N82").Value = -3722
# Row 85
$ws.Range("H85").Value = 2666.6667
$ws.Range("J85").Value = 3000
$ws.Range("L85").Value = 3000
$ws.Range("N85").Value = -5496
# Row 122
$ws.Range("H122").Value = 13640049
$ws.Range("I122").Value = 12502836
$ws.Range("K122").Value = 37508508
$ws.Range("M122").Value = -37506058
# Row 132
$ws.Range("H132").Value = 4853.278
$ws.Range("I132").Value = 4600.091
$ws.Range("K132").Value = 13800.273
$ws.Range("M132").Value = -11270.273
# Row 136
$ws.Range("H136").Value = 24049840
$ws.Range("I136").Value = 37038464
$ws.Range("K136").Value = 111115392
$ws.Range("M136").Value = -111112842

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4834
$ws.Range("I62").Value = 4501
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 4501
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -3877
$ws.Range("N62").Value = -6748
# Row 65
$ws.Range("H65").Value = 4834
$ws.Range("I65").Value = 4501
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 22505
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -19385
$ws.Range("N65").Value = -33740
# Row 70
$ws.Range("H70").Value = 67357.14
$ws.Range("J70").Value = 14300
$ws.Range("L70").Value = 14300
$ws.Range("N70").Value = -14930
# Row 73
$ws.Range("H73").Value = 67357.14
$ws.Range("J73").Value = 14300
$ws.Range("L73").Value = 14300
$ws.Range("N73").Value = -16484
# Row 96
$ws.Range("H96").Value = 1678.7142
$ws.Range("I96").Value = 1716.2
$ws.Range("K96").Value = 1716.2
$ws.Range("M96").Value = -343.2
# Row 132
$ws.Range("H132").Value = 1706.0178
$ws.Range("I132").Value = 1027.3422
$ws.Range("K132").Value = 3082.0266
$ws.Range("M132").Value = -552.0266000000001
